$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 6290
$ws.Range("I58").Value = 966.6667
$ws.Range("J58").Value = 8571.429
$ws.Range("K58").Value = 2900.0001
$ws.Range("L58").Value = 25714.287
$ws.Range("M58").Value = -2750.0001
$ws.Range("N58").Value = -26014.287

$ws.Range("H86").Value = 68190.336
$ws.Range("I86").Value = 100986.6
$ws.Range("J86").Value = 2597.8
$ws.Range("K86").Value = 100986.6
$ws.Range("L86").Value = 2597.8
$ws.Range("M86").Value = -99863.60000000001
$ws.Range("N86").Value = -4843.8

$ws.Range("H87").Value = 25909.092
$ws.Range("I87").Value = 20000
$ws.Range("J87").Value = 26500
$ws.Range("K87").Value = 20000
$ws.Range("L87").Value = 26500
$ws.Range("M87").Value = -18752
$ws.Range("N87").Value = -28996

$ws.Range("H88").Value = 11496195
$ws.Range("I88").Value = 2900.6
$ws.Range("J88").Value = 13890631
$ws.Range("K88").Value = 2900.6
$ws.Range("L88").Value = 13890631
$ws.Range("M88").Value = -2494.6
$ws.Range("N88").Value = -13891443

$ws.Range("H89").Value = 68190.336
$ws.Range("I89").Value = 100986.6
$ws.Range("J89").Value = 2597.8
$ws.Range("K89").Value = 504933
$ws.Range("L89").Value = 12989
$ws.Range("M89").Value = -499317
$ws.Range("N89").Value = -24221

$ws.Range("H90").Value = 25909.092
$ws.Range("I90").Value = 20000
$ws.Range("J90").Value = 26500
$ws.Range("K90").Value = 60000
$ws.Range("L90").Value = 79500
$ws.Range("M90").Value = -53760
$ws.Range("N90").Value = -91980

$ws.Range("H91").Value = 11496195
$ws.Range("I91").Value = 2900.6
$ws.Range("J91").Value = 13890631
$ws.Range("K91").Value = 2900.6
$ws.Range("L91").Value = 13890631
$ws.Range("M91").Value = -1496.6
$ws.Range("N91").Value = -13893439

$ws.Range("H101").Value = 1046.3667
$ws.Range("I101").Value = 519.58826
$ws.Range("J101").Value = 1735.2307
$ws.Range("K101").Value = 1558.76478
$ws.Range("L101").Value = 5205.6921
$ws.Range("M101").Value = 63.23522000000003
$ws.Range("N101").Value = -8449.6921

$ws.Range("H115").Value = 328.33334
$ws.Range("I115").Value = 342.5
$ws.Range("J115").Value = 300
$ws.Range("K115").Value = 1027.5
$ws.Range("L115").Value = 900
$ws.Range("M115").Value = 539.5
$ws.Range("N115").Value = -4034

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1331.6154
$ws.Range("I45").Value = 1435.1666
$ws.Range("J45").Value = 1242.8572
$ws.Range("K45").Value = 1435.1666
$ws.Range("L45").Value = 1242.8572
$ws.Range("M45").Value = -1058.1666
$ws.Range("N45").Value = -1996.8572

$ws.Range("H61").Value = 1785.8695
$ws.Range("I61").Value = 1014.6667
$ws.Range("J61").Value = 2627.182
$ws.Range("K61").Value = 1014.6667
$ws.Range("L61").Value = 2627.182
$ws.Range("M61").Value = -802.6667

$ws.Range("H136").Value = 1785.8695
$ws.Range("I136").Value = 1014.6667
$ws.Range("J136").Value = 2627.182
$ws.Range("K136").Value = 3044.0001
$ws.Range("L136").Value = 7881.545999999999
$ws.Range("M136").Value = -494.0001000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 32631.428
$ws.Range("I82").Value = 3257
$ws.Range("J82").Value = 37527.168
$ws.Range("K82").Value = 3257
$ws.Range("L82").Value = 37527.168
$ws.Range("M82").Value = -2874
$ws.Range("N82").Value = -38293.168

$ws.Range("H85").Value = 32631.428
$ws.Range("I85").Value = 3257
$ws.Range("J85").Value = 37527.168
$ws.Range("K85").Value = 3257
$ws.Range("L85").Value = 37527.168
$ws.Range("M85").Value = -1931
$ws.Range("N85").Value = -40179.168

$ws.Range("H97").Value = 4076
$ws.Range("I97").Value = 3335.5
$ws.Range("J97").Value = 10000
$ws.Range("K97").Value = 3335.5
$ws.Range("L97").Value = 10000
$ws.Range("M97").Value = -2344.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 990
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 975
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 975
$ws.Range("M16").Value = -713
$ws.Range("N16").Value = -1549

$ws.Range("H107").Value = 636
$ws.Range("I107").Value = 590.1539
$ws.Range("J107").Value = 690.1818
$ws.Range("K107").Value = 590.1539
$ws.Range("L107").Value = 690.1818
$ws.Range("M107").Value = 1329.8461
$ws.Range("N107").Value = -4530.1818

$ws.Range("H113").Value = 990
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 975
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 975
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -5315

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 543.8036
$ws.Range("I113").Value = 514.1212
$ws.Range("J113").Value = 586.3913
$ws.Range("K113").Value = 1542.3636
$ws.Range("L113").Value = 1759.1739
$ws.Range("M113").Value = 627.6363999999999

$ws.Range("H131").Value = 5296382.5
$ws.Range("I131").Value = 9529
$ws.Range("J131").Value = 11111921
$ws.Range("K131").Value = 28587
$ws.Range("L131").Value = 33335763
$ws.Range("M131").Value = -23547
$ws.Range("N131").Value = -33345843

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8785.714
$ws.Range("I80").Value = 3500
$ws.Range("J80").Value = 12750
$ws.Range("K80").Value = 3500
$ws.Range("L80").Value = 12750
$ws.Range("M80").Value = -2502
$ws.Range("N80").Value = -14746

$ws.Range("H83").Value = 8785.714
$ws.Range("I83").Value = 3500
$ws.Range("J83").Value = 12750
$ws.Range("K83").Value = 17500
$ws.Range("L83").Value = 63750
$ws.Range("M83").Value = -12508
$ws.Range("N83").Value = -73734

$ws.Range("H97").Value = 3000
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 3000
$ws.Range("N97").Value = -3992

$ws.Range("N113").ClearContents()
$ws.Range("H113").Value = 250000000
$ws.Range("I113").Value = 250000000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 250000000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -249997830

$ws.Range("H126").Value = 3031775.8
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 4168316.5
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 12504949.5
$ws.Range("M126").Value = -530
$ws.Range("N126").Value = -12509889.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1272
$ws.Range("I68").Value = 1253.3334
$ws.Range("J68").Value = 1300
$ws.Range("K68").Value = 1253.3334
$ws.Range("L68").Value = 1300
$ws.Range("M68").Value = -504.3334
$ws.Range("N68").Value = -2798

$ws.Range("H71").Value = 1272
$ws.Range("I71").Value = 1253.3334
$ws.Range("J71").Value = 1300
$ws.Range("K71").Value = 6266.666999999999
$ws.Range("L71").Value = 6500
$ws.Range("M71").Value = -2522.666999999999
$ws.Range("N71").Value = -13988

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2399.4
$ws.Range("I62").Value = 2149.5
$ws.Range("J62").Value = 2566
$ws.Range("K62").Value = 2149.5
$ws.Range("L62").Value = 2566
$ws.Range("M62").Value = -1525.5
$ws.Range("N62").Value = -3814

$ws.Range("H65").Value = 2399.4
$ws.Range("I65").Value = 2149.5
$ws.Range("J65").Value = 2566
$ws.Range("K65").Value = 10747.5
$ws.Range("L65").Value = 12830
$ws.Range("M65").Value = -7627.5
$ws.Range("N65").Value = -19070

Write-Host "Applied Garuda_Profits updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets"
